$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 10-12 ("Viewer" rows) down to rows 11-13 to make room for the
# new "User" row that gets inserted at row 10. Done bottom-up, capturing each
# source row's values first, then pasting the source row's formatting
# (xlPasteFormats = -4122) into the destination row before writing the
# captured values back in (PasteSpecial's xlPasteAll does not reliably carry
# styles in this runtime, and Range.Value is not readable here - Value2 is
# used instead).

# Row 12 -> Row 13
$a12 = $ws.Range("A12").Value2
$b12 = $ws.Range("B12").Value2
$c12 = $ws.Range("C12").Value2
$d12 = $ws.Range("D12").Value2
$ws.Range("A12:D12").Copy()
$ws.Range("A13:D13").PasteSpecial(-4122)
$ws.Range("A13").Value2 = $a12
$ws.Range("B13").Value2 = $b12
$ws.Range("C13").Value2 = $c12
$ws.Range("D13").Value2 = $d12

# Row 11 -> Row 12
$a11 = $ws.Range("A11").Value2
$b11 = $ws.Range("B11").Value2
$c11 = $ws.Range("C11").Value2
$d11 = $ws.Range("D11").Value2
$ws.Range("A11:D11").Copy()
$ws.Range("A12:D12").PasteSpecial(-4122)
$ws.Range("A12").Value2 = $a11
$ws.Range("B12").Value2 = $b11
$ws.Range("C12").Value2 = $c11
$ws.Range("D12").Value2 = $d11

# Row 10 -> Row 11
$a10 = $ws.Range("A10").Value2
$b10 = $ws.Range("B10").Value2
$c10 = $ws.Range("C10").Value2
$d10 = $ws.Range("D10").Value2
$ws.Range("A10:D10").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)
$ws.Range("A11").Value2 = $a10
$ws.Range("B11").Value2 = $b10
$ws.Range("C11").Value2 = $c10
$ws.Range("D11").Value2 = $d10

# Fill the now-vacant row 10 with the new "User" story, reusing the
# formatting of another "User" row (row 9).
$ws.Range("A9:D9").Copy()
$ws.Range("A10:D10").PasteSpecial(-4122)
$ws.Range("A10").Value2 = "User"
$ws.Range("B10").Value2 = "Sort My Capsules"
$ws.Range("C10").Value2 = "Medium"
$ws.Range("D10").Value2 = "To be started"

$ws.Range("G6").Select()
